$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.006228626045210847
$ws.Range("E2").Value = 0.006228626045210847

$ws.Range("D3").Value = 0.9999999999920712
$ws.Range("E3").Value = 0.9999999999920712

$ws.Range("D4").Value = 0.9953874575032584
$ws.Range("E4").Value = 0.9953874575032584

$ws.Range("D5").Value = [double]"2.472840588310253E-41"
$ws.Range("E5").Value = [double]"2.472840588310253E-41"

$ws.Range("D6").Value = 0.9976456668290483
$ws.Range("E6").Value = 0.9976456668290483

$ws.Range("D7").Value = 0.9999999999998348
$ws.Range("E7").Value = [double]"1.652011860642233E-13"

$ws.Range("D8").Value = 0.9922928101930094
$ws.Range("E8").Value = 0.007707189806990589

$ws.Range("D9").Value = 0.9999860414286413
$ws.Range("E9").Value = [double]"1.395857135866319E-05"

$ws.Range("D11").Value = 0.9999999901068849
$ws.Range("E11").Value = [double]"9.893115104908645E-09"
$ws.Range("F11").Value = 3.700501203536987
